$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 held a stray/incorrect figure - the corrected export leaves it blank.
$ws.Range("D3").Value = ""

# Row 7 ("Other") is relabelled "Biogas" with its corrected demand value,
# and a fresh "Other" row (row 8) is appended below it with its own
# corrected value. Formats are cloned from row 7 so the new row picks up
# the same (bold / bordered / centered) label style and blank placeholder
# cells used throughout the table.
$ws.Range("A7:D7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 57.50839259361909

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 24.15205157227206
